# Scheduled-runner market price refresh for Famfrit_Profits sheets
# Updates cached currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# across ALC/ARM/CRP/CUL/GSM/LTW/WVR with freshly fetched Universalis data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 4640.4287
$ws.Range("I19").Value = 2891.5
$ws.Range("J19").Value = 5340
$ws.Range("K19").Value = 2891.5
$ws.Range("L19").Value = 5340
$ws.Range("M19").Value = -2716.5
$ws.Range("N19").Value = -5690

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 4203.4443
$ws.Range("I45").Value = 3466.4
$ws.Range("J45").Value = 5124.75
$ws.Range("K45").Value = 3466.4
$ws.Range("L45").Value = 5124.75
$ws.Range("M45").Value = -3089.4
$ws.Range("N45").Value = -5878.75

# Row 122
$ws.Range("H122").Value = 3200.4285
$ws.Range("I122").Value = 2949.3333
$ws.Range("J122").Value = 3388.75
$ws.Range("K122").Value = 8847.999899999999
$ws.Range("L122").Value = 10166.25
$ws.Range("M122").Value = -6397.999899999999
$ws.Range("N122").Value = -15066.25

# Row 129
$ws.Range("H129").Value = 74999.5
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 74999.5
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 74999.5
$ws.Range("N129").Value = -84999.5

$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 409
$ws.Range("I23").Value = 409
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 409
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -169

# Row 27
$ws.Range("H27").Value = 409
$ws.Range("I27").Value = 409
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 409
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -217

# Row 141
$ws.Range("H141").Value = 83946.5
$ws.Range("I141").Value = 27648
$ws.Range("J141").Value = 90201.89
$ws.Range("K141").Value = 27648
$ws.Range("L141").Value = 90201.89
$ws.Range("M141").Value = -22468
$ws.Range("N141").Value = -100561.89

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 1168.5333
$ws.Range("I113").Value = 468.63635
$ws.Range("J113").Value = 3093.25
$ws.Range("K113").Value = 1405.90905
$ws.Range("L113").Value = 9279.75
$ws.Range("M113").Value = 764.09095
$ws.Range("N113").Value = -13619.75

$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 6586.2
$ws.Range("I21").Value = 10411
$ws.Range("J21").Value = 849
$ws.Range("K21").Value = 10411
$ws.Range("L21").Value = 849
$ws.Range("M21").Value = -10238
$ws.Range("N21").Value = -1195

# Row 29
$ws.Range("H29").Value = 6135.6665
$ws.Range("I29").Value = 1000
$ws.Range("J29").Value = 8703.5
$ws.Range("K29").Value = 1000
$ws.Range("L29").Value = 8703.5
$ws.Range("M29").Value = -710
$ws.Range("N29").Value = -9283.5

# Row 30
$ws.Range("H30").Value = 6586.2
$ws.Range("I30").Value = 10411
$ws.Range("J30").Value = 849
$ws.Range("K30").Value = 10411
$ws.Range("L30").Value = 849
$ws.Range("M30").Value = -10306
$ws.Range("N30").Value = -1059

# Row 102
$ws.Range("H102").Value = 4970.125
$ws.Range("I102").Value = 2643.0715
$ws.Range("J102").Value = 8228
$ws.Range("K102").Value = 2643.0715
$ws.Range("L102").Value = 8228
$ws.Range("M102").Value = -1021.0715
$ws.Range("N102").Value = -11472

# Row 122
$ws.Range("H122").Value = 1829.2858
$ws.Range("I122").Value = 1112.7059
$ws.Range("J122").Value = 4874.75
$ws.Range("K122").Value = 3338.1177
$ws.Range("L122").Value = 14624.25
$ws.Range("M122").Value = -888.1176999999998
$ws.Range("N122").Value = -19524.25

$ws = $wb.Worksheets.Item("LTW")
# Row 5
$ws.Range("H5").Value = 30000
$ws.Range("I5").Value = 30000
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 30000
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -29887
$ws.Range("N5").ClearContents()

# Row 7
$ws.Range("H7").Value = 7562.125
$ws.Range("I7").Value = 6998.6
$ws.Range("J7").Value = 8501.333000000001
$ws.Range("K7").Value = 6998.6
$ws.Range("L7").Value = 8501.333000000001
$ws.Range("M7").Value = -6886.6
$ws.Range("N7").Value = -8725.333000000001

# Row 36
$ws.Range("H36").Value = 50000
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 50000
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 50000
$ws.Range("N36").Value = -51124

# Row 40
$ws.Range("H40").Value = 6911.2856
$ws.Range("I40").Value = 4600
$ws.Range("J40").Value = 9993
$ws.Range("K40").Value = 4600
$ws.Range("L40").Value = 9993
$ws.Range("M40").Value = -4464
$ws.Range("N40").Value = -10265

# Row 122
$ws.Range("H122").Value = 4130.943
$ws.Range("I122").Value = 3783.92
$ws.Range("J122").Value = 4998.5
$ws.Range("K122").Value = 11351.76
$ws.Range("L122").Value = 14995.5
$ws.Range("M122").Value = -8901.76
$ws.Range("N122").Value = -19895.5

# Row 126
$ws.Range("H126").Value = 7562.125
$ws.Range("I126").Value = 6998.6
$ws.Range("J126").Value = 8501.333000000001
$ws.Range("K126").Value = 20995.8
$ws.Range("L126").Value = 25503.999
$ws.Range("M126").Value = -18525.8
$ws.Range("N126").Value = -30443.999

# Row 130
$ws.Range("H130").Value = 61496.25
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 61496.25
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 61496.25
$ws.Range("N130").Value = -71536.25

$ws = $wb.Worksheets.Item("WVR")
# Row 18
$ws.Range("H18").Value = 204738
$ws.Range("I18").Value = 1000947
$ws.Range("J18").Value = 5685.75
$ws.Range("K18").Value = 1000947
$ws.Range("L18").Value = 5685.75
$ws.Range("M18").Value = -1000774
$ws.Range("N18").Value = -6031.75

# Row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()

# Row 24
$ws.Range("H24").Value = 19000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 19000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 19000
$ws.Range("N24").Value = -19460

# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").ClearContents()

# Row 81
$ws.Range("H81").Value = 3595.3572
$ws.Range("I81").Value = 3195.25
$ws.Range("J81").Value = 5996
$ws.Range("K81").Value = 6390.5
$ws.Range("L81").Value = 11992
$ws.Range("M81").Value = -5329.5
$ws.Range("N81").Value = -14114

# Row 84
$ws.Range("H84").Value = 3595.3572
$ws.Range("I84").Value = 3195.25
$ws.Range("J84").Value = 5996
$ws.Range("K84").Value = 31952.5
$ws.Range("L84").Value = 59960
$ws.Range("M84").Value = -26648.5
$ws.Range("N84").Value = -70568

# Row 122
$ws.Range("H122").Value = 3646.8147
$ws.Range("I122").Value = 2693.762
$ws.Range("J122").Value = 6982.5
$ws.Range("K122").Value = 8081.286
$ws.Range("L122").Value = 20947.5
$ws.Range("M122").Value = -5631.286
$ws.Range("N122").Value = -25847.5

# Row 133
$ws.Range("H133").Value = 79666.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 79666.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 79666.5
$ws.Range("N133").Value = -89786.5
